# Natmi following Dr Hou advice
# Rebuild the LR-pairs data table to include the new "ECs" cell type,
# expanding the sending/target cluster cross-product from 2x2 to 3x3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A="ECs";  D="ECs";  E=2; F=0.6666666666666666; G=1.853892;          H=5.561676;          I=0.5711238486747862; J=0.571123848674786;  K=1; L=0.3333333333333333; M=0.1604616666666667; N=0.481385;  O=0.01863590595952956; P=0.01863590595952956; Q=0.29747860014;      R=2.67730740126;      S=0.0106434103351479;  T=0.0106434103351479 },
    @{ A="ECs";  D="FAPs"; E=2; F=0.6666666666666666; G=1.853892;          H=5.561676;          I=0.5711238486747862; J=0.571123848674786;  K=3; L=1;                  M=7.100470666666666; N=21.301412; O=0.8246437068815905;  P=0.8246437068815904;  Q=13.163505765168;    R=118.471551886512;   S=0.4709736876596562;  T=0.470973687659656  },
    @{ A="ECs";  D="sCs";  E=2; F=0.6666666666666666; G=1.853892;          H=5.561676;          I=0.5711238486747862; J=0.571123848674786;  K=3; L=1;                  M=1.349417333333333; N=4.048252000000001; O=0.15672038715888;   P=0.15672038715888;    Q=2.501673998928001;  R=22.515065990352;    S=0.08950675067998211; T=0.08950675067998207 },
    @{ A="FAPs"; D="ECs";  E=3; F=1;                  G=0.4262446666666667; H=1.278734;          I=0.1313121230922664; J=0.1313121230922664; K=1; L=0.3333333333333333; M=0.1604616666666667; N=0.481385;  O=0.01863590595952956; P=0.01863590595952956; Q=0.06839592962111112; R=0.61556336659;      S=0.002447120377293647; T=0.002447120377293646 },
    @{ A="FAPs"; D="FAPs"; E=3; F=1;                  G=0.4262446666666667; H=1.278734;          I=0.1313121230922664; J=0.1313121230922664; K=3; L=1;                  M=7.100470666666666; N=21.301412; O=0.8246437068815905;  P=0.8246437068815904;  Q=3.026537752489777;  R=27.238839772408;    S=0.1082857159452983;  T=0.1082857159452982  },
    @{ A="FAPs"; D="sCs";  E=3; F=1;                  G=0.4262446666666667; H=1.278734;          I=0.1313121230922664; J=0.1313121230922664; K=3; L=1;                  M=1.349417333333333; N=4.048252000000001; O=0.15672038715888;   P=0.15672038715888;    Q=0.575181941440889;  R=5.176637472968001;  S=0.0205792867696745;  T=0.0205792867696745  },
    @{ A="sCs";  D="ECs";  E=3; F=1;                  G=0.9659053333333333; H=2.897716;          I=0.2975640282329475; J=0.2975640282329475; K=1; L=0.3333333333333333; M=0.1604616666666667; N=0.481385;  O=0.01863590595952956; P=0.01863590595952956; Q=0.1549907796288889;  R=1.39491701666;      S=0.005545375247088007; T=0.005545375247088007 },
    @{ A="sCs";  D="FAPs"; E=3; F=1;                  G=0.9659053333333333; H=2.897716;          I=0.2975640282329475; J=0.2975640282329475; K=3; L=1;                  M=7.100470666666666; N=21.301412; O=0.8246437068815905;  P=0.8246437068815904;  Q=6.858382486110221;  R=61.725442374992;    S=0.245384303276636;  T=0.245384303276636  },
    @{ A="sCs";  D="sCs";  E=3; F=1;                  G=0.9659053333333333; H=2.897716;          I=0.2975640282329475; J=0.2975640282329475; K=3; L=1;                  M=1.349417333333333; N=4.048252000000001; O=0.15672038715888;   P=0.15672038715888;    Q=1.303409399159111;  R=11.730684592432;    S=0.04663434970922344; T=0.04663434970922343 }
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $row = $startRow + $i

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = "Dhh"
    $ws.Cells.Item($row, 3).Value = "Hhip"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
